$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The Cases (row 2) and Samples (row 3) tabs pointed at the stale
# BCellLymphoma / Lymphoma input files; repoint them at the OsteoSarcoma
# Neo4j/Web data files used by the rest of this test-case workbook.
$ws.Range("D2").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_WebData.xlsx"
$ws.Range("D3").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC05_Canine_Filter_Diagnosis-OsteoSarcoma_WebData.xlsx"

# Leave the sheet with the same selection/scroll state saved in the file.
$ws.Activate()
$ws.Range("E4").Select()
